$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1: "Our Stand" (Heading1) -> Title-styled paragraph, text split
# into three runs: "Our" / " " / "Stand" (mirrors pandoc's per-word runs).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Style = "Title"

$rOur = $p1.Range.Duplicate
$rOur.Find.Execute("Our", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$rSpace1 = $p1.Range.Duplicate
$rSpace1.Find.Execute(" ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# Re-assert the (unchanged) font size on just the space so it gets written out
# as its own run, splitting "Our" / " " / "Stand" into separate <w:r> runs.
$rSpace1.Font.Size = $rSpace1.Font.Size

# ---------------------------------------------------------------------------
# Paragraph 2: "By Dorothy Day" (bold, Normal) -> "Dorothy Day" (Authors
# style, not bold), text split into three runs: "Dorothy" / " " / "Day".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("By Dorothy Day", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dorothy Day", 2) | Out-Null

$p2 = $d.Paragraphs.Item(2)
$p2.Style = "Authors"

$rSpace2 = $p2.Range.Duplicate
$rSpace2.Find.Execute(" ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rSpace2.Font.Size = $rSpace2.Font.Size

Write-Host "p1:" $p1.Range.Text "| style=" $p1.Style.NameLocal
Write-Host "p2:" $p2.Range.Text "| style=" $p2.Style.NameLocal
